$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -7.399799999999999
$ws.Range("D7").Value = -7.666200000000003
$ws.Range("D16").Value = -8.704200000000005
$ws.Range("D28").Value = -8.283199999999999
$ws.Range("D29").Value = -7.240100000000001
$ws.Range("D32").Value = -9.042399999999999
$ws.Range("D40").Value = -8.185899999999991
$ws.Range("D52").Value = -7.646399999999997
$ws.Range("D57").Value = -8.5284
$ws.Range("D66").Value = -7.2276
$ws.Range("D100").Value = -8.358200000000005
